$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{ rowNum=2; A="ECs"; B="Timp3"; C="Agtr2"; D="ECs"; E=3; F=1; G=62.15419033333333; H=186.462571; I=0.5307382952913039; J=0.5798280707535227; K=1; L=0.3333333333333333; M=0.01476566666666667; N=0.044297; O=0.003369932031170864; P=0.003537680284170835; Q=0.9177480563985556; R=8.259732507587; S=0.001788551981471185; T=0.002051246334113549 }
    @{ rowNum=3; A="ECs"; B="Timp3"; C="Agtr2"; D="FAPs"; E=3; F=1; G=62.15419033333333; H=186.462571; I=0.5307382952913039; J=0.5798280707535227; K=3; L=1; M=3.743532666666667; N=11.230598; O=0.8543773151546027; P=0.8969064524470826; Q=232.6762418830509; R=2094.086176947458; S=0.4534507597807149; T=0.520051537968778 }
    @{ rowNum=4; A="FAPs"; B="Timp3"; C="Agtr2"; D="MuSCs"; E=3; F=1; G=62.15419033333333; H=186.462571; I=0.5307382952913039; J=0.5798280707535227; K=2; L=1; M=0.6232935; N=1.246587; O=0.1422527528142264; P=0.09955586726874661; Q=38.74030283252949; R=232.441816995177; S=0.07549898352911778; T=0.05772528645063113 }
    @{ rowNum=5; A="FAPs"; B="Timp3"; C="Agtr2"; D="ECs"; E=3; F=1; G=25.210481; H=75.631443; I=0.2152737834352902; J=0.2351851819258409; K=1; L=0.3333333333333333; M=0.01476566666666667; N=0.044297; O=0.003369932031170864; P=0.003537680284170835; Q=0.3722495589523334; R=3.350246030571; S=0.0007254580182699241; T=0.0008320099812281784 }
    @{ rowNum=6; A="FAPs"; B="Timp3"; C="Agtr2"; D="FAPs"; E=3; F=1; G=25.210481; H=75.631443; I=0.2152737834352902; J=0.2351851819258409; K=3; L=1; M=3.743532666666667; N=11.230598; O=0.8543773151546027; P=0.8969064524470826; Q=94.37625916587935; R=849.3863324929141; S=0.1839250371146166; T=0.2109391071892277 }
    @{ rowNum=7; A="FAPs"; B="Timp3"; C="Agtr2"; D="MuSCs"; E=3; F=1; G=25.210481; H=75.631443; I=0.2152737834352902; J=0.2351851819258409; K=2; L=1; M=0.6232935; N=1.246587; O=0.1422527528142264; P=0.09955586726874661; Q=15.7135289391735; R=94.281173635041; S=0.03062328830240365; T=0.02341406475538504 }
    @{ rowNum=8; A="MuSCs"; B="Timp3"; C="Agtr2"; D="ECs"; E=2; F=1; G=29.744252; H=59.488504; I=0.2539879212734059; J=0.1849867473206364; K=1; L=0.3333333333333333; M=0.01476566666666667; N=0.044297; O=0.003369932031170864; P=0.003537680284170835; Q=0.4391937102813334; R=2.635162261688; S=0.0008559220314297542; T=0.0006544239688291074 }
    @{ rowNum=9; A="MuSCs"; B="Timp3"; C="Agtr2"; D="FAPs"; E=2; F=1; G=29.744252; H=59.488504; I=0.2539879212734059; J=0.1849867473206364; K=3; L=1; M=3.743532666666667; N=11.230598; O=0.8543773151546027; P=0.8969064524470826; Q=111.3485790075653; R=668.091474045392; S=0.2170015182592711; T=0.1659158072890768 }
    @{ rowNum=10; A="MuSCs"; B="Timp3"; C="Agtr2"; D="MuSCs"; E=2; F=1; G=29.744252; H=59.488504; I=0.2539879212734059; J=0.1849867473206364; K=2; L=1; M=0.6232935; N=1.246587; O=0.1422527528142264; P=0.09955586726874661; Q=18.539398933962; R=74.15759573584799; S=0.03613048098270502; T=0.01841651606273044 }
)

foreach ($row in $rowsData) {
    $ws.Range("A" + $row.rowNum).Value = $row.A
    $ws.Range("B" + $row.rowNum).Value = $row.B
    $ws.Range("C" + $row.rowNum).Value = $row.C
    $ws.Range("D" + $row.rowNum).Value = $row.D
    $ws.Range("E" + $row.rowNum).Value = $row.E
    $ws.Range("F" + $row.rowNum).Value = $row.F
    $ws.Range("G" + $row.rowNum).Value = $row.G
    $ws.Range("H" + $row.rowNum).Value = $row.H
    $ws.Range("I" + $row.rowNum).Value = $row.I
    $ws.Range("J" + $row.rowNum).Value = $row.J
    $ws.Range("K" + $row.rowNum).Value = $row.K
    $ws.Range("L" + $row.rowNum).Value = $row.L
    $ws.Range("M" + $row.rowNum).Value = $row.M
    $ws.Range("N" + $row.rowNum).Value = $row.N
    $ws.Range("O" + $row.rowNum).Value = $row.O
    $ws.Range("P" + $row.rowNum).Value = $row.P
    $ws.Range("Q" + $row.rowNum).Value = $row.Q
    $ws.Range("R" + $row.rowNum).Value = $row.R
    $ws.Range("S" + $row.rowNum).Value = $row.S
    $ws.Range("T" + $row.rowNum).Value = $row.T
}